$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4851.6665
$ws.Range("I64").Value = 3795
$ws.Range("J64").Value = 5908.3335
$ws.Range("K64").Value = 3795
$ws.Range("L64").Value = 5908.3335
$ws.Range("M64").Value = -3547
$ws.Range("N64").Value = -6404.3335
$ws.Range("H67").Value = 4851.6665
$ws.Range("I67").Value = 3795
$ws.Range("J67").Value = 5908.3335
$ws.Range("K67").Value = 3795
$ws.Range("L67").Value = 5908.3335
$ws.Range("M67").Value = -2937
$ws.Range("N67").Value = -7624.3335
$ws.Range("H112").Value = 772645.3
$ws.Range("J112").Value = 836949.0600000001
$ws.Range("L112").Value = 2510847.18
$ws.Range("N112").Value = -2513063.18
$ws.Range("H138").Value = 3653.1938
$ws.Range("I138").Value = 1476.7778
$ws.Range("J138").Value = 4480.845
$ws.Range("K138").Value = 4430.3334
$ws.Range("L138").Value = 13442.535
$ws.Range("M138").Value = 709.6665999999996
$ws.Range("N138").Value = -23722.535
$ws.Range("H141").Value = 2483.6667
$ws.Range("I141").Value = 1693
$ws.Range("K141").Value = 5079
$ws.Range("M141").Value = 101

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5921.89
$ws.Range("I32").Value = 5314.7085
$ws.Range("K32").Value = 5314.7085
$ws.Range("M32").Value = -5027.7085
$ws.Range("H45").Value = 1818.3334
$ws.Range("I45").Value = 1592.24
$ws.Range("J45").Value = 2524.875
$ws.Range("K45").Value = 1592.24
$ws.Range("L45").Value = 2524.875
$ws.Range("M45").Value = -1215.24
$ws.Range("N45").Value = -3278.875
$ws.Range("H74").Value = 2699.3447
$ws.Range("I74").Value = 1742.3684
$ws.Range("K74").Value = 1742.3684
$ws.Range("M74").Value = -868.3684000000001
$ws.Range("H75").Value = 99078.5
$ws.Range("J75").Value = 98000
$ws.Range("L75").Value = 98000
$ws.Range("N75").Value = -99748
$ws.Range("H77").Value = 2699.3447
$ws.Range("I77").Value = 1742.3684
$ws.Range("K77").Value = 8711.842000000001
$ws.Range("M77").Value = -4343.842000000001
$ws.Range("H78").Value = 99078.5
$ws.Range("J78").Value = 98000
$ws.Range("L78").Value = 294000
$ws.Range("N78").Value = -302736
$ws.Range("H97").Value = 451.85184
$ws.Range("I97").Value = 451.85184
$ws.Range("K97").Value = 451.85184
$ws.Range("M97").Value = 44.14816000000002
$ws.Range("H102").Value = 8004104.5
$ws.Range("I102").Value = 4060.6
$ws.Range("J102").Value = 40004280
$ws.Range("K102").Value = 4060.6
$ws.Range("L102").Value = 40004280
$ws.Range("M102").Value = -2438.6
$ws.Range("N102").Value = -40007524
$ws.Range("H122").Value = 11159.454
$ws.Range("I122").Value = 13187.286
$ws.Range("K122").Value = 39561.858
$ws.Range("M122").Value = -37111.858

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 5190.3335
$ws.Range("I37").Value = 1287.75
$ws.Range("J37").Value = 12995.5
$ws.Range("K37").Value = 1287.75
$ws.Range("L37").Value = 12995.5
$ws.Range("M37").Value = -1150.75
$ws.Range("N37").Value = -13269.5
$ws.Range("H86").Value = 4451465
$ws.Range("I86").Value = 7416940
$ws.Range("J86").Value = 3251.8333
$ws.Range("K86").Value = 7416940
$ws.Range("L86").Value = 3251.8333
$ws.Range("M86").Value = -7415817
$ws.Range("N86").Value = -5497.8333
$ws.Range("H89").Value = 4451465
$ws.Range("I89").Value = 7416940
$ws.Range("J89").Value = 3251.8333
$ws.Range("K89").Value = 37084700
$ws.Range("L89").Value = 16259.1665
$ws.Range("M89").Value = -37079084
$ws.Range("N89").Value = -27491.1665
$ws.Range("H94").Value = 1062.3636
$ws.Range("I94").Value = 1085
$ws.Range("J94").Value = 898.25
$ws.Range("K94").Value = 1085
$ws.Range("L94").Value = 898.25
$ws.Range("M94").Value = -634
$ws.Range("N94").Value = -1800.25
$ws.Range("H99").Value = 2952.7917
$ws.Range("I99").Value = 2726.9285
$ws.Range("J99").Value = 3269
$ws.Range("K99").Value = 2726.9285
$ws.Range("L99").Value = 3269
$ws.Range("M99").Value = -1228.9285
$ws.Range("N99").Value = -6265
$ws.Range("H105").Value = 5931.385
$ws.Range("I105").Value = 5360.9
$ws.Range("K105").Value = 5360.9
$ws.Range("M105").Value = -3613.9
$ws.Range("H134").Value = 4723.2324
$ws.Range("I134").Value = 4974.758
$ws.Range("K134").Value = 14924.274
$ws.Range("M134").Value = -12389.274

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2211.3572
$ws.Range("I31").Value = 1861.8529
$ws.Range("K31").Value = 1861.8529
$ws.Range("M31").Value = -1566.8529
$ws.Range("H34").Value = 2211.3572
$ws.Range("I34").Value = 1861.8529
$ws.Range("K34").Value = 1861.8529
$ws.Range("M34").Value = -1659.8529
$ws.Range("H62").Value = 4023.12
$ws.Range("I62").Value = 3176.4546
$ws.Range("K62").Value = 3176.4546
$ws.Range("M62").Value = -2552.4546
$ws.Range("H65").Value = 4023.12
$ws.Range("I65").Value = 3176.4546
$ws.Range("K65").Value = 15882.273
$ws.Range("M65").Value = -12762.273

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 6608.6665
$ws.Range("I56").Value = 6608.6665
$ws.Range("K56").Value = 6608.6665
$ws.Range("M56").Value = -6078.6665
$ws.Range("H129").Value = 84848
$ws.Range("I129").Value = 100937.7
$ws.Range("K129").Value = 302813.1
$ws.Range("M129").Value = -297813.1
$ws.Range("H132").Value = 2997.9211
$ws.Range("J132").Value = 3208.375
$ws.Range("L132").Value = 28875.375
$ws.Range("N132").Value = -33935.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5166.3335
$ws.Range("I80").Value = 5166.3335
$ws.Range("K80").Value = 5166.3335
$ws.Range("M80").Value = -4168.3335
$ws.Range("H83").Value = 5166.3335
$ws.Range("I83").Value = 5166.3335
$ws.Range("K83").Value = 25831.6675
$ws.Range("M83").Value = -20839.6675
$ws.Range("H100").Value = 62669
$ws.Range("J100").Value = 50118
$ws.Range("L100").Value = 50118
$ws.Range("N100").Value = -52282
$ws.Range("H113").Value = 2459.6
$ws.Range("I113").Value = 2133
$ws.Range("K113").Value = 2133
$ws.Range("M113").Value = 37
$ws.Range("H122").Value = 4971.3335
$ws.Range("I122").Value = 4320.0713
$ws.Range("J122").Value = 7250.75
$ws.Range("K122").Value = 12960.2139
$ws.Range("L122").Value = 21752.25
$ws.Range("M122").Value = -10510.2139
$ws.Range("N122").Value = -26652.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 5300
$ws.Range("J3").Value = 7450
$ws.Range("L3").Value = 7450
$ws.Range("N3").Value = -7674
$ws.Range("H9").Value = 374.66666
$ws.Range("I9").Value = 292.25
$ws.Range("J9").Value = 539.5
$ws.Range("K9").Value = 292.25
$ws.Range("L9").Value = 539.5
$ws.Range("M9").Value = -68.25
$ws.Range("N9").Value = -987.5
$ws.Range("H15").Value = 5300
$ws.Range("J15").Value = 7450
$ws.Range("L15").Value = 7450
$ws.Range("N15").Value = -7790
$ws.Range("H20").Value = 679664.7
$ws.Range("J20").Value = 679664.7
$ws.Range("L20").Value = 679664.7
$ws.Range("N20").Value = -680116.7
$ws.Range("H22").Value = 2847.2144
$ws.Range("J22").Value = 3151.4167
$ws.Range("L22").Value = 3151.4167
$ws.Range("N22").Value = -3741.4167
$ws.Range("H27").Value = 2847.2144
$ws.Range("J27").Value = 3151.4167
$ws.Range("L27").Value = 3151.4167
$ws.Range("N27").Value = -3365.4167
$ws.Range("H55").Value = 533.875
$ws.Range("I55").Value = 516.5
$ws.Range("K55").Value = 516.5
$ws.Range("M55").Value = -343.5
$ws.Range("H122").Value = 7840.1177
$ws.Range("I122").Value = 7339.7
$ws.Range("J122").Value = 8555
$ws.Range("K122").Value = 22019.1
$ws.Range("L122").Value = 25665
$ws.Range("M122").Value = -19569.1
$ws.Range("N122").Value = -30565
$ws.Range("H132").Value = 3809.1191
$ws.Range("I132").Value = 3529.0435
$ws.Range("K132").Value = 10587.1305
$ws.Range("M132").Value = -8057.130500000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 36550.863
$ws.Range("I62").Value = 42825.74
$ws.Range("J62").Value = 12497.167
$ws.Range("K62").Value = 42825.74
$ws.Range("L62").Value = 12497.167
$ws.Range("M62").Value = -42201.74
$ws.Range("N62").Value = -13745.167
$ws.Range("H65").Value = 36550.863
$ws.Range("I65").Value = 42825.74
$ws.Range("J65").Value = 12497.167
$ws.Range("K65").Value = 214128.7
$ws.Range("L65").Value = 62485.835
$ws.Range("M65").Value = -211008.7
$ws.Range("N65").Value = -68725.83499999999
$ws.Range("H96").Value = 3397.5
$ws.Range("I96").Value = 4175
$ws.Range("K96").Value = 4175
$ws.Range("M96").Value = -2802
$ws.Range("H132").Value = 2365.3076
$ws.Range("I132").Value = 2181.2354
$ws.Range("J132").Value = 2713
$ws.Range("K132").Value = 6543.706200000001
$ws.Range("L132").Value = 8139
$ws.Range("M132").Value = -4013.706200000001
$ws.Range("N132").Value = -13199
